$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contribuyente")

$ws.Cells.Item(2, 9).Value = "DK7331645124473461205164"
$ws.Cells.Item(2, 10).Value = "JMD00@Agua2024.com"
$ws.Cells.Item(3, 9).Value = "ES8265614874165615445616"
$ws.Cells.Item(3, 10).Value = "RFM00@Agua2024.com"
$ws.Cells.Item(4, 9).Value = "RO8832569523016220165156"
$ws.Cells.Item(4, 10).Value = "ALM00@Agua2024.com"
$ws.Cells.Item(5, 9).Value = "DE7424561937521546497521"
$ws.Cells.Item(5, 10).Value = "ALN00@Agua2024.com"
$ws.Cells.Item(6, 9).Value = "MC6436520125638451012515"
$ws.Cells.Item(6, 10).Value = "RGO00@Agua2024.com"
$ws.Cells.Item(7, 9).Value = "ES0721584976902154655487"
$ws.Cells.Item(7, 10).Value = "SDO00@Agua2024.com"
$ws.Cells.Item(8, 9).Value = "GR9420125003305201112544"
$ws.Cells.Item(8, 10).Value = "VDO00@Agua2024.com"
$ws.Cells.Item(9, 9).Value = "ES2821651484690980008984"
$ws.Cells.Item(9, 10).Value = "VBP00@Agua2024.com"
$ws.Cells.Item(10, 9).Value = "FI5620960043043554600000"
$ws.Cells.Item(10, 10).Value = "ABP00@Agua2024.com"
$ws.Cells.Item(11, 9).Value = "ES7921564975243245467995"
$ws.Cells.Item(11, 10).Value = "MAP00@Agua2024.com"
$ws.Cells.Item(12, 9).Value = "LT8032566221522587754554"
$ws.Cells.Item(12, 10).Value = "ACP00@Agua2024.com"
$ws.Cells.Item(13, 9).Value = "EE2023215465315456411515"
$ws.Cells.Item(13, 10).Value = "LBP00@Agua2024.com"
$ws.Cells.Item(16, 9).Value = "SM2125894363475485700145"
$ws.Cells.Item(16, 10).Value = "LBR00@Agua2024.com"
$ws.Cells.Item(17, 9).Value = "ES9596431245118150005156"
$ws.Cells.Item(17, 10).Value = "SBR00@Agua2024.com"
$ws.Cells.Item(18, 9).Value = "AT6825030000114574745458"
$ws.Cells.Item(18, 10).Value = "AGR00@Agua2024.com"
$ws.Cells.Item(19, 9).Value = "IT8915953684811254695203"
$ws.Cells.Item(19, 10).Value = "DGR00@Agua2024.com"
$ws.Cells.Item(21, 9).Value = "DK5800750184310702510000"
$ws.Cells.Item(21, 10).Value = "CIS00@Agua2024.com"
$ws.Cells.Item(22, 9).Value = "ES5023455254943263234457"
$ws.Cells.Item(22, 10).Value = "MQ00@Agua2024.com"
$ws.Cells.Item(23, 9).Value = "GR4920910936583000000000"
$ws.Cells.Item(23, 10).Value = "BV00@Agua2024.com"
$ws.Cells.Item(24, 9).Value = "ES3720960043032159000000"
$ws.Cells.Item(24, 10).Value = "PP00@Agua2024.com"
$ws.Cells.Item(25, 9).Value = "DE5512669681115112121210"
$ws.Cells.Item(25, 10).Value = "PC00@Agua2024.com"
$ws.Cells.Item(27, 9).Value = "ES2956187775315550000651"
$ws.Cells.Item(27, 10).Value = "GMM00@Agua2024.com"
$ws.Cells.Item(28, 9).Value = "ES0425516848021156151054"
$ws.Cells.Item(28, 10).Value = "CGM00@Agua2024.com"
$ws.Cells.Item(29, 9).Value = "PT5764578946740051516490"
$ws.Cells.Item(29, 10).Value = "GMM01@Agua2024.com"
$ws.Cells.Item(30, 9).Value = "ES4534698752714600549403"
$ws.Cells.Item(30, 10).Value = "CSN00@Agua2024.com"
$ws.Cells.Item(31, 9).Value = "ES2766649444162310000255"
$ws.Cells.Item(31, 10).Value = "ALO00@Agua2024.com"
$ws.Cells.Item(32, 9).Value = "FR5623185484465641685100"
$ws.Cells.Item(32, 10).Value = "GMM02@Agua2024.com"
$ws.Cells.Item(36, 9).Value = "DE5021508149175421346497"
$ws.Cells.Item(36, 10).Value = "TCP00@Agua2024.com"
$ws.Cells.Item(37, 9).Value = "DE6721346154503164978451"
$ws.Cells.Item(37, 10).Value = "CCP00@Agua2024.com"
$ws.Cells.Item(38, 9).Value = "ES7225187786311225455548"
$ws.Cells.Item(38, 10).Value = "CAP00@Agua2024.com"
$ws.Cells.Item(40, 9).Value = "ES2396536214865214585214"
$ws.Cells.Item(40, 10).Value = "ELR00@Agua2024.com"
$ws.Cells.Item(41, 9).Value = "ES6885461325251978750005"
$ws.Cells.Item(41, 10).Value = "DLR00@Agua2024.com"
$ws.Cells.Item(43, 9).Value = "ES5020960043073071400000"
$ws.Cells.Item(43, 10).Value = "SPR00@Agua2024.com"
$ws.Cells.Item(44, 9).Value = "ES8220960043042158800000"
$ws.Cells.Item(44, 10).Value = "GPR00@Agua2024.com"
$ws.Cells.Item(45, 9).Value = "ES7521654587985156484454"
$ws.Cells.Item(45, 10).Value = "SAS00@Agua2024.com"
$ws.Cells.Item(46, 9).Value = "ES3251651681961210656510"
$ws.Cells.Item(46, 10).Value = "BFS00@Agua2024.com"
$ws.Cells.Item(47, 9).Value = "ES5566552211148855332200"
$ws.Cells.Item(47, 10).Value = "DFG00@Agua2024.com"
$ws.Cells.Item(48, 9).Value = "GB9720910936583000000000"
$ws.Cells.Item(48, 10).Value = "GMG00@Agua2024.com"
$ws.Cells.Item(49, 9).Value = "DE9301821135910205540000"
$ws.Cells.Item(49, 10).Value = "MLG00@Agua2024.com"
$ws.Cells.Item(50, 9).Value = "DE7822631245526916432102"
$ws.Cells.Item(50, 10).Value = "IAG00@Agua2024.com"
$ws.Cells.Item(51, 9).Value = "ES2120960043043075700000"
$ws.Cells.Item(51, 10).Value = "IAG01@Agua2024.com"
$ws.Cells.Item(52, 9).Value = "SM7325635478321002541225"
$ws.Cells.Item(52, 10).Value = "GPG00@Agua2024.com"
$ws.Cells.Item(53, 9).Value = "ES6832154697195423121000"
$ws.Cells.Item(53, 10).Value = "RGG00@Agua2024.com"
$ws.Cells.Item(55, 9).Value = "GB5520008521528775113366"
$ws.Cells.Item(55, 10).Value = "AFG00@Agua2024.com"
$ws.Cells.Item(60, 9).Value = "ES8020960043033000100000"
$ws.Cells.Item(60, 10).Value = "APM00@Agua2024.com"
$ws.Cells.Item(61, 9).Value = "GB0836585214290025478551"
$ws.Cells.Item(61, 10).Value = "GMM03@Agua2024.com"
$ws.Cells.Item(62, 9).Value = "ES9012548523465214585214"
$ws.Cells.Item(62, 10).Value = "ALM01@Agua2024.com"
$ws.Cells.Item(63, 9).Value = "ES6931624561042546920007"
$ws.Cells.Item(63, 10).Value = "LDM00@Agua2024.com"
$ws.Cells.Item(64, 9).Value = "ES1436154231712500312566"
$ws.Cells.Item(64, 10).Value = "AD00@Agua2024.com"
$ws.Cells.Item(65, 9).Value = "ES8244875664127231645789"
$ws.Cells.Item(65, 10).Value = "LR00@Agua2024.com"
$ws.Cells.Item(66, 9).Value = "ES7920960031442124800000"
$ws.Cells.Item(66, 10).Value = "LC00@Agua2024.com"
$ws.Cells.Item(68, 9).Value = "ES1933218885441445121022"
$ws.Cells.Item(68, 10).Value = "AFG01@Agua2024.com"
$ws.Cells.Item(69, 9).Value = "ES8462581542713690044508"
$ws.Cells.Item(69, 10).Value = "DGG00@Agua2024.com"
$ws.Cells.Item(71, 9).Value = "ES8020960043033000100000"
$ws.Cells.Item(71, 10).Value = "SOG00@Agua2024.com"
$ws.Cells.Item(72, 9).Value = "PT3536952365020014425254"
$ws.Cells.Item(72, 10).Value = "VVG00@Agua2024.com"
$ws.Cells.Item(73, 9).Value = "ES9565168874641561561500"
$ws.Cells.Item(73, 10).Value = "VMG00@Agua2024.com"
$ws.Cells.Item(74, 9).Value = "ES3220960583831234500000"
$ws.Cells.Item(74, 10).Value = "MBG00@Agua2024.com"
$ws.Cells.Item(75, 9).Value = "ES7221416325811510005514"
$ws.Cells.Item(75, 10).Value = "CBG00@Agua2024.com"
$ws.Cells.Item(79, 9).Value = "ES8163516541828944000984"
$ws.Cells.Item(79, 10).Value = "SDM00@Agua2024.com"
$ws.Cells.Item(80, 8).Value = "23658965274585223202"
$ws.Cells.Item(80, 9).Value = "ES6223658965274585223202"
$ws.Cells.Item(80, 10).Value = "EGM00@Agua2024.com"
$ws.Cells.Item(81, 9).Value = "FI6132658012367712548745"
$ws.Cells.Item(81, 10).Value = "GPM00@Agua2024.com"
$ws.Cells.Item(82, 9).Value = "ES7223652365142254222000"
$ws.Cells.Item(82, 10).Value = "EAM00@Agua2024.com"
$ws.Cells.Item(84, 9).Value = "ES9232584216971684051000"
$ws.Cells.Item(84, 10).Value = "MG00@Agua2024.com"
$ws.Cells.Item(86, 9).Value = "ES7395485212315484010000"
$ws.Cells.Item(86, 10).Value = "JAC00@Agua2024.com"
$ws.Cells.Item(87, 9).Value = "LT9321856333126985542360"
$ws.Cells.Item(87, 10).Value = "BDC00@Agua2024.com"
$ws.Cells.Item(88, 9).Value = "ES5736245978133245679001"
$ws.Cells.Item(88, 10).Value = "NGC00@Agua2024.com"
$ws.Cells.Item(89, 9).Value = "ES7631245164156597845124"
$ws.Cells.Item(89, 10).Value = "MLC00@Agua2024.com"
$ws.Cells.Item(91, 9).Value = "SE6832574512085411002255"
$ws.Cells.Item(91, 10).Value = "MLC01@Agua2024.com"
$ws.Cells.Item(92, 9).Value = "ES4420960043013468900000"
$ws.Cells.Item(92, 10).Value = "CFC00@Agua2024.com"
$ws.Cells.Item(93, 9).Value = "ES5631215643855060225021"
$ws.Cells.Item(93, 10).Value = "CGC00@Agua2024.com"
$ws.Cells.Item(95, 9).Value = "ES1665165654918886005001"
$ws.Cells.Item(95, 10).Value = "KSC00@Agua2024.com"
$ws.Cells.Item(102, 8).Value = "65645150865168448896"
$ws.Cells.Item(102, 9).Value = "AT8365645150865168448896"
$ws.Cells.Item(102, 10).Value = "MHC00@Agua2024.com"
$ws.Cells.Item(103, 8).Value = "26551681807651415636"
$ws.Cells.Item(103, 9).Value = "IT3526551681807651415636"
$ws.Cells.Item(103, 10).Value = "CLD00@Agua2024.com"
$ws.Cells.Item(104, 9).Value = "HU2399558741836555551120"
$ws.Cells.Item(104, 10).Value = "MFD00@Agua2024.com"
$ws.Cells.Item(106, 8).Value = "51556584221251000254"
$ws.Cells.Item(106, 9).Value = "IE6851556584221251000254"
$ws.Cells.Item(106, 10).Value = "GMM04@Agua2024.com"
$ws.Cells.Item(128, 8).Value = "62541122421110105611"
$ws.Cells.Item(128, 9).Value = "LT9362541122421110105611"
$ws.Cells.Item(128, 10).Value = "DMC00@Agua2024.com"
$ws.Cells.Item(129, 9).Value = "ES6855065688761051056105"
$ws.Cells.Item(129, 10).Value = "EBC00@Agua2024.com"
$ws.Cells.Item(131, 9).Value = "ES9712548521518742146695"
$ws.Cells.Item(131, 10).Value = "MSC00@Agua2024.com"
$ws.Cells.Item(132, 9).Value = "ES9001826530120201560000"
$ws.Cells.Item(132, 10).Value = "MDC00@Agua2024.com"
$ws.Cells.Item(133, 9).Value = "ES9021651651812511133551"
$ws.Cells.Item(133, 10).Value = "MFC00@Agua2024.com"
$ws.Cells.Item(134, 9).Value = "ES6851651487910005118185"
$ws.Cells.Item(134, 10).Value = "CDD00@Agua2024.com"
$ws.Cells.Item(136, 9).Value = "AT3122515651915640081000"
$ws.Cells.Item(136, 10).Value = "HPD00@Agua2024.com"
